$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values scraped on 2023-01-28 (GitHub Actions symbol-list refresh).
# Columns B (Coin) and C (Link) are plain text; D (Price) and E (Volume 1h) look
# numeric/percentage so we force text format first to keep them as strings like
# the source data (avoids Excel auto-converting "307.95" / "0.93%" into numbers).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '307.95'
Set-TextValue $ws.Range("E2") '0.93%'
Set-TextValue $ws.Range("E3") '8.16%'
Set-TextValue $ws.Range("D4") '5.099'
Set-TextValue $ws.Range("E4") '1.47%'
Set-TextValue $ws.Range("D5") '0.08116'
Set-TextValue $ws.Range("E5") '1.34%'
Set-TextValue $ws.Range("D6") '1.969'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D7") '4.207'
Set-TextValue $ws.Range("E7") '1.51%'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws.Range("D8") '7.950'
Set-TextValue $ws.Range("E8") '2.17%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D9") '0.9291'
Set-TextValue $ws.Range("E9") '1.00%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D10") '0.1426'
Set-TextValue $ws.Range("E10") '12.43%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D11") '0.1961'
Set-TextValue $ws.Range("E11") '2.53%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D12") '0.09019'
Set-TextValue $ws.Range("E12") '-1.29%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D13") '0.03510'
Set-TextValue $ws.Range("E13") '1.44%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D14") '0.09830'
Set-TextValue $ws.Range("E14") '-0.30%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D15") '0.001413'
Set-TextValue $ws.Range("E15") '0.75%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D16") '0.006107'
Set-TextValue $ws.Range("E16") '-2.18%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D17") '3.681'
Set-TextValue $ws.Range("E17") '-4.35%'
Set-TextValue $ws.Range("D18") '3.475'
Set-TextValue $ws.Range("E18") '2.83%'
Set-TextValue $ws.Range("E19") '1.28%'
Set-TextValue $ws.Range("D20") '0.1303'
Set-TextValue $ws.Range("E20") '-3.31%'
Set-TextValue $ws.Range("D21") '4.817'
Set-TextValue $ws.Range("E21") '-7.53%'
Set-TextValue $ws.Range("D22") '0.2455'
Set-TextValue $ws.Range("E22") '6.47%'
Set-TextValue $ws.Range("D23") '0.04421'
Set-TextValue $ws.Range("D24") '0.001218'
Set-TextValue $ws.Range("E24") '-1.34%'
Set-TextValue $ws.Range("E25") '-1.16%'
Set-TextValue $ws.Range("E27") '4.06%'
Set-TextValue $ws.Range("D39") '0.02096'
Set-TextValue $ws.Range("E39") '8.18%'
Set-TextValue $ws.Range("D40") '0.05126'
Set-TextValue $ws.Range("E40") '-1.82%'
Set-TextValue $ws.Range("D41") '0.007471'
Set-TextValue $ws.Range("E41") '-1.90%'
Set-TextValue $ws.Range("D42") '0.01014'
Set-TextValue $ws.Range("E42") '-0.16%'
Set-TextValue $ws.Range("E43") '0.82%'
Set-TextValue $ws.Range("D44") '0.002134'
Set-TextValue $ws.Range("E44") '-0.88%'
Set-TextValue $ws.Range("D45") '0.009250'
Set-TextValue $ws.Range("E45") '-3.87%'
Set-TextValue $ws.Range("D46") '0.00006254'
Set-TextValue $ws.Range("E46") '0.76%'
Set-TextValue $ws.Range("E47") '0.02%'
Set-TextValue $ws.Range("D48") '0.003061'
Set-TextValue $ws.Range("D49") '0.001602'
Set-TextValue $ws.Range("E49") '-3.55%'
Set-TextValue $ws.Range("E50") '0.02%'
Set-TextValue $ws.Range("E51") '0.02%'
